$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap column widths: column A takes the old column B width, column B takes the old column A width.
# ColumnWidth snaps to this runtimes pixel grid, so these inputs are chosen to land on the
# closest achievable stored width to the target (16.42578125 / 15.7109375).
$ws.Columns.Item(1).ColumnWidth = 15.65
$ws.Columns.Item(2).ColumnWidth = 14.85

# Update cell values (A1:B32)
$ws.Range("A1").Value = -0.14506580613478093
$ws.Range("B1").Value = 0.14498677983044672
$ws.Range("A2").Value = -0.13464194998036572
$ws.Range("B2").Value = 0.13435507665639168
$ws.Range("A3").Value = -0.084652088605764675
$ws.Range("B3").Value = 0.084474637856807533
$ws.Range("A4").Value = -0.076474637904423659
$ws.Range("B4").Value = 0.075947045498379495
$ws.Range("A5").Value = -0.072947045523830134
$ws.Range("B5").Value = 0.07114452309267616
$ws.Range("A6").Value = -0.024260166790220694
$ws.Range("B6").Value = 0.023820441661269953
$ws.Range("A7").Value = -0.014728806993261845
$ws.Range("B7").Value = 0.014624639725566269
$ws.Range("A8").Value = -0.0046246397926661409
$ws.Range("B8").Value = 0.0044434993410376045
$ws.Range("A9").Value = -0.0024434993694124607
$ws.Range("B9").Value = 0.0022969989683838676
$ws.Range("A10").Value = -0.00029699899704382915
$ws.Range("B10").Value = 0.000288188548399404
$ws.Range("A11").Value = 0.0027118114180026964
$ws.Range("B11").Value = -0.0027271454879187473
$ws.Range("A12").Value = 0.006227145451947802
$ws.Range("B12").Value = -0.006340913928725378
$ws.Range("A13").Value = 0.0098409138935249274
$ws.Range("B13").Value = -0.0098934578944067653
$ws.Range("A14").Value = 0.017893457837587334
$ws.Range("B14").Value = -0.017918040062753882
$ws.Range("A15").Value = 0.018918040041169704
$ws.Range("B15").Value = -0.018937692006016249
$ws.Range("A16").Value = -0.0060341884282912872
$ws.Range("B16").Value = 0.0060034659909899979
$ws.Range("A17").Value = -0.00400346601700452
$ws.Range("B17").Value = 0.0039999999639555028
$ws.Range("A18").Value = -0.0043448298803134833
$ws.Range("B18").Value = 0.00431007367116365
$ws.Range("A19").Value = -0.00031007369158331244
$ws.Range("B19").Value = 0.000075708580593580876
$ws.Range("A20").Value = 0.0039242913987447281
$ws.Range("B20").Value = -0.0039807140451983969
$ws.Range("A21").Value = -0.0040056864334800935
$ws.Range("B21").Value = 0.0039999999792286189
$ws.Range("A22").Value = -0.045702988076255835
$ws.Range("B22").Value = 0.045492374363776733
$ws.Range("A23").Value = -0.040492374395949327
$ws.Range("B23").Value = 0.040097710870688807
$ws.Range("A24").Value = -0.020097710979539052
$ws.Range("B24").Value = 0.019999999889711795
$ws.Range("A25").Value = -0.04488435632886123
$ws.Range("B25").Value = 0.044832677165343782
$ws.Range("A26").Value = -0.042332677194904633
$ws.Range("B26").Value = 0.04226944840971214
$ws.Range("A27").Value = -0.039769448440234889
$ws.Range("B27").Value = 0.039415978840048105
$ws.Range("A28").Value = -0.057725300274372948
$ws.Range("B28").Value = 0.057300468947713235
$ws.Range("A29").Value = -0.050300469008615956
$ws.Range("B29").Value = 0.050186413268680674
$ws.Range("A30").Value = 0.0098135864082360413
$ws.Range("B30").Value = -0.0099431114459616055
$ws.Range("A31").Value = -0.014022950871735773
$ws.Range("B31").Value = 0.014001302050660414
$ws.Range("A32").Value = -0.004001302127990769
$ws.Range("B32").Value = 0.0039999999520095031
